{"js": "// Replace the header date and each \"NNN\u00d7N=\" expression in the practice\n// table with its updated value. All old strings are unique within the\n// document, so a simple exact-text search+replace per pair is sufficient\n// and keeps the original run formatting (rFonts/sz) intact.\nconst replacements = [\n  [\"2025-12-24 Wednesday\", \"2025-12-25 Thursday\"],\n  [\"322\u00d74=\", \"109\u00d74=\"],\n  [\"997\u00d75=\", \"740\u00d79=\"],\n  [\"590\u00d79=\", \"881\u00d76=\"],\n  [\"544\u00d75=\", \"204\u00d74=\"],\n  [\"529\u00d79=\", \"905\u00d76=\"],\n  [\"952\u00d78=\", \"519\u00d77=\"],\n  [\"403\u00d77=\", \"135\u00d73=\"],\n  [\"955\u00d74=\", \"639\u00d77=\"],\n  [\"102\u00d77=\", \"382\u00d76=\"],\n  [\"909\u00d73=\", \"871\u00d78=\"],\n  [\"606\u00d76=\", \"122\u00d73=\"],\n  [\"168\u00d78=\", \"431\u00d77=\"],\n  [\"654\u00d73=\", \"586\u00d77=\"],\n  [\"613\u00d77=\", \"158\u00d73=\"],\n  [\"826\u00d75=\", \"607\u00d77=\"],\n  [\"744\u00d79=\", \"152\u00d78=\"],\n  [\"452\u00d73=\", \"465\u00d76=\"],\n  [\"212\u00d76=\", \"570\u00d79=\"],\n  [\"108\u00d77=\", \"598\u00d75=\"],\n  [\"162\u00d79=\", \"582\u00d77=\"],\n  [\"723\u00d74=\", \"154\u00d74=\"],\n  [\"695\u00d77=\", \"752\u00d78=\"],\n  [\"806\u00d76=\", \"981\u00d77=\"],\n  [\"659\u00d72=\", \"276\u00d78=\"],\n  [\"530\u00d72=\", \"818\u00d72=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the header date and each \"NNN\u00d7N=\" expression in the practice\n# table with its updated value. Every old string is unique within the\n# document, so Find/Execute with Replace:=wdReplaceOne (2) on the whole\n# document Content range is sufficient and preserves the existing run\n# formatting (rFonts/sz) of the text it replaces.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-12-24 Wednesday\", \"2025-12-25 Thursday\"),\n    @(\"322\u00d74=\", \"109\u00d74=\"),\n    @(\"997\u00d75=\", \"740\u00d79=\"),\n    @(\"590\u00d79=\", \"881\u00d76=\"),\n    @(\"544\u00d75=\", \"204\u00d74=\"),\n    @(\"529\u00d79=\", \"905\u00d76=\"),\n    @(\"952\u00d78=\", \"519\u00d77=\"),\n    @(\"403\u00d77=\", \"135\u00d73=\"),\n    @(\"955\u00d74=\", \"639\u00d77=\"),\n    @(\"102\u00d77=\", \"382\u00d76=\"),\n    @(\"909\u00d73=\", \"871\u00d78=\"),\n    @(\"606\u00d76=\", \"122\u00d73=\"),\n    @(\"168\u00d78=\", \"431\u00d77=\"),\n    @(\"654\u00d73=\", \"586\u00d77=\"),\n    @(\"613\u00d77=\", \"158\u00d73=\"),\n    @(\"826\u00d75=\", \"607\u00d77=\"),\n    @(\"744\u00d79=\", \"152\u00d78=\"),\n    @(\"452\u00d73=\", \"465\u00d76=\"),\n    @(\"212\u00d76=\", \"570\u00d79=\"),\n    @(\"108\u00d77=\", \"598\u00d75=\"),\n    @(\"162\u00d79=\", \"582\u00d77=\"),\n    @(\"723\u00d74=\", \"154\u00d74=\"),\n    @(\"695\u00d77=\", \"752\u00d78=\"),\n    @(\"806\u00d76=\", \"981\u00d77=\"),\n    @(\"659\u00d72=\", \"276\u00d78=\"),\n    @(\"530\u00d72=\", \"818\u00d72=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $r = $d.Content\n    $r.Find.ClearFormatting()\n    $r.Find.Replacement.ClearFormatting()\n    $found = $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        Write-Host (\"WARNING: text not found: \" + $oldText)\n    }\n}\n\n$d.Saved = $false\n"}
